{"js": "// Office.js (Word JavaScript API) edit script.\n// Replaces the answer text in specific table cells (two-digit division\n// problems) with newly generated problems, matching the positions used\n// in the source diff. Each cell is addressed by (row, col) in the single\n// table that holds the worksheet, and we verify the existing text before\n// replacing it so the edit only touches the intended run.\n\nconst edits = [{\"row\": 0, \"col\": 0, \"old\": \"32\u00f74=8, 0\", \"new\": \"75\u00f74=18, 3\"}, {\"row\": 0, \"col\": 1, \"old\": \"95\u00f77=13, 4\", \"new\": \"39\u00f77=5, 4\"}, {\"row\": 0, \"col\": 2, \"old\": \"58\u00f79=6, 4\", \"new\": \"49\u00f73=16, 1\"}, {\"row\": 0, \"col\": 3, \"old\": \"95\u00f73=31, 2\", \"new\": \"10\u00f77=1, 3\"}, {\"row\": 0, \"col\": 4, \"old\": \"45\u00f76=7, 3\", \"new\": \"17\u00f78=2, 1\"}, {\"row\": 4, \"col\": 0, \"old\": \"78\u00f79=8, 6\", \"new\": \"97\u00f74=24, 1\"}, {\"row\": 4, \"col\": 1, \"old\": \"31\u00f79=3, 4\", \"new\": \"60\u00f72=30, 0\"}, {\"row\": 4, \"col\": 2, \"old\": \"72\u00f75=14, 2\", \"new\": \"24\u00f74=6, 0\"}, {\"row\": 4, \"col\": 3, \"old\": \"40\u00f73=13, 1\", \"new\": \"84\u00f75=16, 4\"}, {\"row\": 4, \"col\": 4, \"old\": \"34\u00f79=3, 7\", \"new\": \"64\u00f73=21, 1\"}, {\"row\": 8, \"col\": 0, \"old\": \"50\u00f75=10, 0\", \"new\": \"98\u00f75=19, 3\"}, {\"row\": 8, \"col\": 1, \"old\": \"21\u00f74=5, 1\", \"new\": \"48\u00f73=16, 0\"}, {\"row\": 8, \"col\": 2, \"old\": \"49\u00f74=12, 1\", \"new\": \"27\u00f75=5, 2\"}, {\"row\": 8, \"col\": 3, \"old\": \"15\u00f75=3, 0\", \"new\": \"50\u00f77=7, 1\"}, {\"row\": 8, \"col\": 4, \"old\": \"79\u00f78=9, 7\", \"new\": \"69\u00f73=23, 0\"}, {\"row\": 12, \"col\": 0, \"old\": \"60\u00f74=15, 0\", \"new\": \"38\u00f76=6, 2\"}, {\"row\": 12, \"col\": 1, \"old\": \"75\u00f76=12, 3\", \"new\": \"26\u00f77=3, 5\"}, {\"row\": 12, \"col\": 2, \"old\": \"36\u00f74=9, 0\", \"new\": \"24\u00f72=12, 0\"}, {\"row\": 12, \"col\": 3, \"old\": \"96\u00f72=48, 0\", \"new\": \"65\u00f76=10, 5\"}, {\"row\": 12, \"col\": 4, \"old\": \"85\u00f78=10, 5\", \"new\": \"30\u00f77=4, 2\"}, {\"row\": 16, \"col\": 0, \"old\": \"43\u00f72=21, 1\", \"new\": \"51\u00f75=10, 1\"}, {\"row\": 16, \"col\": 1, \"old\": \"80\u00f77=11, 3\", \"new\": \"44\u00f75=8, 4\"}, {\"row\": 16, \"col\": 2, \"old\": \"22\u00f78=2, 6\", \"new\": \"44\u00f78=5, 4\"}, {\"row\": 16, \"col\": 3, \"old\": \"56\u00f72=28, 0\", \"new\": \"71\u00f76=11, 5\"}, {\"row\": 16, \"col\": 4, \"old\": \"38\u00f76=6, 2\", \"new\": \"36\u00f77=5, 1\"}];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const edit of edits) {\n  const cell = table.getCell(edit.row, edit.col);\n  cell.body.paragraphs.load(\"items\");\n  await context.sync();\n\n  const paragraph = cell.body.paragraphs.items[0];\n  paragraph.load(\"text\");\n  await context.sync();\n\n  if (paragraph.text !== edit.old) {\n    throw new Error(\n      `Unexpected cell text at row ${edit.row}, col ${edit.col}: ` +\n      `expected \"${edit.old}\" but found \"${paragraph.text}\"`\n    );\n  }\n\n  const range = paragraph.getRange();\n  range.insertText(edit.new, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Replaces the answer text in specific table cells (two-digit division\n# problems) with newly generated problems, matching the positions used\n# in the source diff. Each cell is addressed by (Row, Col) (1-indexed,\n# as Word COM expects) in the single table that holds the worksheet, and\n# the existing text is verified before replacing it so the edit only\n# touches the intended run and its formatting (fonts/size/alignment) is\n# left untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$edits = @(\n    @{ Row = 1; Col = 1; Old = '32\u00f74=8, 0'; New = '75\u00f74=18, 3' },\n    @{ Row = 1; Col = 2; Old = '95\u00f77=13, 4'; New = '39\u00f77=5, 4' },\n    @{ Row = 1; Col = 3; Old = '58\u00f79=6, 4'; New = '49\u00f73=16, 1' },\n    @{ Row = 1; Col = 4; Old = '95\u00f73=31, 2'; New = '10\u00f77=1, 3' },\n    @{ Row = 1; Col = 5; Old = '45\u00f76=7, 3'; New = '17\u00f78=2, 1' },\n    @{ Row = 5; Col = 1; Old = '78\u00f79=8, 6'; New = '97\u00f74=24, 1' },\n    @{ Row = 5; Col = 2; Old = '31\u00f79=3, 4'; New = '60\u00f72=30, 0' },\n    @{ Row = 5; Col = 3; Old = '72\u00f75=14, 2'; New = '24\u00f74=6, 0' },\n    @{ Row = 5; Col = 4; Old = '40\u00f73=13, 1'; New = '84\u00f75=16, 4' },\n    @{ Row = 5; Col = 5; Old = '34\u00f79=3, 7'; New = '64\u00f73=21, 1' },\n    @{ Row = 9; Col = 1; Old = '50\u00f75=10, 0'; New = '98\u00f75=19, 3' },\n    @{ Row = 9; Col = 2; Old = '21\u00f74=5, 1'; New = '48\u00f73=16, 0' },\n    @{ Row = 9; Col = 3; Old = '49\u00f74=12, 1'; New = '27\u00f75=5, 2' },\n    @{ Row = 9; Col = 4; Old = '15\u00f75=3, 0'; New = '50\u00f77=7, 1' },\n    @{ Row = 9; Col = 5; Old = '79\u00f78=9, 7'; New = '69\u00f73=23, 0' },\n    @{ Row = 13; Col = 1; Old = '60\u00f74=15, 0'; New = '38\u00f76=6, 2' },\n    @{ Row = 13; Col = 2; Old = '75\u00f76=12, 3'; New = '26\u00f77=3, 5' },\n    @{ Row = 13; Col = 3; Old = '36\u00f74=9, 0'; New = '24\u00f72=12, 0' },\n    @{ Row = 13; Col = 4; Old = '96\u00f72=48, 0'; New = '65\u00f76=10, 5' },\n    @{ Row = 13; Col = 5; Old = '85\u00f78=10, 5'; New = '30\u00f77=4, 2' },\n    @{ Row = 17; Col = 1; Old = '43\u00f72=21, 1'; New = '51\u00f75=10, 1' },\n    @{ Row = 17; Col = 2; Old = '80\u00f77=11, 3'; New = '44\u00f75=8, 4' },\n    @{ Row = 17; Col = 3; Old = '22\u00f78=2, 6'; New = '44\u00f78=5, 4' },\n    @{ Row = 17; Col = 4; Old = '56\u00f72=28, 0'; New = '71\u00f76=11, 5' },\n    @{ Row = 17; Col = 5; Old = '38\u00f76=6, 2'; New = '36\u00f77=5, 1' }\n)\n\nforeach ($edit in $edits) {\n    $cell = $t.Cell($edit.Row, $edit.Col)\n    $r = $cell.Range\n    # Drop the trailing end-of-cell marker so only the visible text is\n    # addressed (keeps paragraph/run formatting intact).\n    $r.MoveEnd(1, -1) | Out-Null\n\n    $current = $r.Text\n    if ($current -ne $edit.Old) {\n        throw \"Unexpected cell text at row $($edit.Row), col $($edit.Col): expected '$($edit.Old)' but found '$current'\"\n    }\n\n    $r.Text = $edit.New\n}\n"}
